$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The old label "A,4PC=1PC" used in the "A" (pieces-per-set) rule column for
# the roll-journal weld-repair rows (44-57) is renamed to "A,1SET=4PC".
$ws.Range("E44:E57").Value = "A,1SET=4PC"

# Leave the active selection where the author left it after scrolling the
# sheet down to this block of rows.
$ws.Range("E61").Select()
